# Effort_3_2021 workbook update — "Successful Effort excel creation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level housekeeping -------------------------------------------------
# Drop the (empty) workbookProtection element.
$wb.Unprotect()

# Restore the normal window geometry that Excel stamps into bookViews
# when a workbook is opened/edited interactively (instead of the
# freshly-created, size-less default).
$win = $excel.ActiveWindow
$win.Left   = 360
$win.Top    = 525
$win.Width  = 19815
$win.Height = 7365

# --- Row 2: fix up the first data row ------------------------------------------
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = "rikushwa"
$ws.Range("D2").Value = "Meetings / Communication"
$ws.Range("E2").Value = "Mail Communication"
$ws.Range("F2").Value = "2021-03-02 22:52:51"
$ws.Range("G2").Value = 1.5
$ws.Range("H2").Value = "Medium"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "TOPSI"

# --- New rows 3-8: additional logged efforts ------------------------------------
$rows = @(
    @{ r=3; a="'0";            d="Service-Task"; e="DSTUM";    g=1.5  },
    @{ r=4; a="SCTASK0073278"; d="Incident";      e="Incident"; g=0.75 },
    @{ r=5; a="INC0597119";    d="Incident";      e="Incident"; g=0.75 },
    @{ r=6; a="INC0580213";    d="Incident";      e="Incident"; g=0.75 },
    @{ r=7; a="SCTASK0073460"; d="Incident";      e="Incident"; g=0.75 },
    @{ r=8; a="INC0597552";    d="Incident";      e="Incident"; g=0.75 }
)

foreach ($row in $rows) {
    $n = $row.r
    $ws.Range("A$n").Value = $row.a
    $ws.Range("B$n").Value = ""
    $ws.Range("C$n").Value = "rikushwa"
    $ws.Range("D$n").Value = $row.d
    $ws.Range("E$n").Value = $row.e
    $ws.Range("F$n").Value = "2021-03-02 22:52:51"
    $ws.Range("G$n").Value = $row.g
    $ws.Range("H$n").Value = "Medium"
    $ws.Range("I$n").Value = ""
    $ws.Range("J$n").Value = ""
    $ws.Range("K$n").Value = "TOPSI"
}

# --- Row 9: grand-total formula --------------------------------------------------
$ws.Range("A9").Formula = "=SUM(G1:G9)"

# --- Header-row fill: re-apply the same colour so it round-trips with a clean alpha channel
$ws.Range("A1:K1").Interior.Color = 3733478

# --- Selection / active cell, matching the saved UI state -----------------------
$ws.Range("M15").Select()
